# Refresh the crypto price/volume table (cryptos.xlsx) with the latest
# scrape, as produced by the "Updated cryptos list" GitHub Actions job.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.008.76'
$ws.Range('E2').Value = '  +0.60%  '
$ws.Range('D3').Value = '1.641.86'
$ws.Range('E3').Value = '  +0.67%  '
$ws.Range('E4').Value = '  +0.35%  '
$ws.Range('D5').Value = "'215.89"
$ws.Range('E5').Value = '  +0.83%  '
$ws.Range('E6').Value = '  +0.18%  '
$ws.Range('E7').Value = '  +0.35%  '
$ws.Range('E8').Value = '  +0.50%  '
$ws.Range('E9').Value = '  +1.03%  '
$ws.Range('E10').Value = '  +0.22%  '
$ws.Range('E11').Value = '  +0.43%  '
$ws.Range('E12').Value = '  +0.71%  '
$ws.Range('E13').Value = '  +0.66%  '
$ws.Range('D14').Value = '1.654.19'
$ws.Range('E14').Value = '  +1.75%  '
$ws.Range('E15').Value = '  +0.15%  '
$ws.Range('D16').Value = '0.0₃0762'
$ws.Range('D18').Value = '26.105.05'
$ws.Range('E18').Value = '  +0.95%  '
$ws.Range('E19').Value = '  +0.37%  '
$ws.Range('D20').Value = "'194.34"
$ws.Range('E20').Value = '  +0.57%  '
$ws.Range('E21').Value = '  -0.59%  '
$ws.Range('E22').Value = '  +0.16%  '
$ws.Range('D23').Value = "'6.20"
$ws.Range('E23').Value = '  -0.96%  '
$ws.Range('E24').Value = '  -0.73%  '
$ws.Range('E25').Value = '  +4.81%  '
$ws.Range('D26').Value = "'1.00"
$ws.Range('E26').Value = '  +0.11%  '
$ws.Range('D27').Value = "'143.03"
$ws.Range('E27').Value = '  -0.29%  '
$ws.Range('E28').Value = '  +0.72%  '
$ws.Range('D29').Value = "'15.53"
$ws.Range('E29').Value = '  +0.84%  '
$ws.Range('E30').Value = '  +0.67%  '
$ws.Range('E31').Value = '  -0.55%  '
$ws.Range('E32').Value = '  -0.13%  '
$ws.Range('D33').Value = "'3.26"
$ws.Range('E33').Value = '  +1.23%  '
$ws.Range('D34').Value = "'1.54"
$ws.Range('E34').Value = '  -0.78%  '
$ws.Range('E35').Value = '  +1.19%  '
$ws.Range('E36').Value = '  +0.52%  '
$ws.Range('D37').Value = '1.130.45'
$ws.Range('E37').Value = '  -0.65%  '
$ws.Range('E38').Value = '  -0.90%  '
$ws.Range('E39').Value = '  -0.57%  '
$ws.Range('D40').Value = "'0.0157"
$ws.Range('E40').Value = '  +0.39%  '
$ws.Range('E41').Value = '  +0.92%  '
$ws.Range('D42').Value = "'99.14"
$ws.Range('E42').Value = '  -0.08%  '
$ws.Range('E43').Value = '  +0.10%  '
$ws.Range('D44').Value = '1.778.41'
$ws.Range('E44').Value = '  +0.73%  '
$ws.Range('E45').Value = '  +4.98%  '
$ws.Range('D46').Value = "'56.50"
$ws.Range('E46').Value = '  +0.57%  '
$ws.Range('B47').Value = 'RenderToken'
$ws.Range('C47').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D47').Value = "'1.49"
$ws.Range('E47').Value = '  +3.66%  '
$ws.Range('B48').Value = 'Cronos'
$ws.Range('C48').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D48').Value = "'0.0522"
$ws.Range('E48').Value = '  -1.25%  '
$ws.Range('D49').Value = "'7.75"
$ws.Range('E49').Value = '  +1.78%  '
$ws.Range('E50').Value = '  -0.21%  '
$ws.Range('E51').Value = '  +0.28%  '
